# Add a new column F containing a repeated image-URL string in rows 1-14,
# mirroring the hyperlink-style "BARCODE generator" source columns already
# present in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

$url = "https://t4.ftcdn.net/jpg/01/33/48/03/360_F_133480376_PWlsZ1Bdr2SVnTRpb8jCtY59CyEBdoUt.jpg"

for ($r = 1; $r -le 14; $r++) {
    $ws.Cells.Item($r, 6).Value = $url
}

# Widen the new column to fit the long URL text.
$ws.Range("F1:F14").ColumnWidth = 81.05338541666667

# Leave the new range selected, as it was when the author saved the file.
$ws.Range("F1:F14").Select()
